$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------------
# 1) Refresh the cached "datetimeFigureOut" field text (slide master + every
#    custom layout) from 8/15/2018 -> 8/29/2018.
# ---------------------------------------------------------------------------
$design = $p.Designs.Item(1)
$master = $design.SlideMaster

for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $shp = $master.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        if ($shp.TextFrame.TextRange.Text -eq "8/15/2018") {
            $shp.TextFrame.TextRange.Text = "8/29/2018"
        }
    }
}

for ($l = 1; $l -le $master.CustomLayouts.Count; $l++) {
    $layout = $master.CustomLayouts.Item($l)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $shp = $layout.Shapes.Item($i)
        if ($shp.HasTextFrame) {
            if ($shp.TextFrame.TextRange.Text -eq "8/15/2018") {
                $shp.TextFrame.TextRange.Text = "8/29/2018"
            }
        }
    }
}

# ---------------------------------------------------------------------------
# 2) Slide shape text / geometry touch-ups.
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)

    if ($shp.Name -eq "TextBox 3") {
        # "In-memory buffer" -> "In-memory Buffer"
        $shp.TextFrame.TextRange.Text = "In-memory Buffer"
    }

    if ($shp.Name -eq "Rectangle 109") {
        # "Message processing" -> "Message Processing", plus a small
        # reflow of the textbox's position/size.
        $shp.TextFrame.TextRange.Text = "Message Processing"
        $shp.Left = 496.358031496063
        $shp.Width = 178.64275590551182
    }
}
